$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# ---------------------------------------------------------------------------
# 1) Insert the new "ExceptionsFolder" row right before the existing
#    "CompletedFolder" row (old worksheet row 43). This pushes everything
#    from old row 43 onward down by one.
# ---------------------------------------------------------------------------
$ws.Rows.Item(43).Insert()

$ws.Range("A43").Value = "ExceptionsFolder"
$ws.Range("B43").Value = "\\cavmfil001\Common\SinglePaymentScheme\Cross Compliance {0}\9) Processing Folders\Robot\{1}\Exceptions"
# The row copied formatting (including a C column style) from the row above;
# the new row only has Name/Value, so drop the stray formatted C cell.
$ws.Range("C43").Clear()

# ---------------------------------------------------------------------------
# 2) Insert the new "outlookExceptionsFolder" row right after the existing
#    "inProgressFolder" row. After step 1, "inProgressFolder" now sits at
#    worksheet row 53, so the new row goes in at row 54.
# ---------------------------------------------------------------------------
$ws.Rows.Item(54).Insert()

$ws.Range("A54").Value = "outlookExceptionsFolder"
$ws.Range("B54").Value = "Exceptions"

# ---------------------------------------------------------------------------
# 3) Grow Table1 so it covers the two newly-inserted rows (A1:C82 -> A1:C84).
# ---------------------------------------------------------------------------
$lo.Resize($ws.Range("A1:C84"))

# ---------------------------------------------------------------------------
# 4) Update the view: scroll/select to roughly match where the author left
#    the cursor after the edit (row 43 and row 54 highlighted).
# ---------------------------------------------------------------------------
[void]$ws.Rows.Item(54).Select()
[void]$ws.Rows.Item(43).Select()
